# Rotate the species-record data among rows 5, 6 and 7:
#   old row 6 -> row 5
#   old row 7 -> row 6
#   old row 5 -> row 7
# Only columns A, B, E, F, G, H, Q, R change; everything else in these
# rows is identical between the three rows already.
# Note: this runtime's Range.Value getter is unreliable, so Value2 is
# used for both reading and writing cell contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

# Capture the original values for rows 5, 6, 7 before overwriting anything.
$orig5 = @{}
$orig6 = @{}
$orig7 = @{}
foreach ($col in $cols) {
    $orig5[$col] = $ws.Range($col + "5").Value2
    $orig6[$col] = $ws.Range($col + "6").Value2
    $orig7[$col] = $ws.Range($col + "7").Value2
}

# Apply the cyclic rotation: row6 -> row5, row7 -> row6, row5 -> row7
foreach ($col in $cols) {
    $ws.Range($col + "5").Value2 = $orig6[$col]
    $ws.Range($col + "6").Value2 = $orig7[$col]
    $ws.Range($col + "7").Value2 = $orig5[$col]
}
